$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.010643243789673
$ws.Range("B1").Value = 1.01332700252533
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.608077168464661
$ws.Range("E1").Value = 0.9706447124481201
